$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.001214305068890198
$ws.Range("D2").Value = 0.002046462813289907
$ws.Range("E2").Value = 0.4304973541418349
$ws.Range("F2").Value = 0.8133486523056774
$ws.Range("G2").Value = 0.002341840673284616
$ws.Range("I2").Value = 0.6112125958755144
$ws.Range("O2").Value = 2.719297177951887

$ws.Range("C3").Value = 0.001082874213185647
$ws.Range("D3").Value = 0.001855303547699805
$ws.Range("E3").Value = 0.3750918814099862
$ws.Range("F3").Value = 0.7525647962886808
$ws.Range("G3").Value = 0.002346712352702368
$ws.Range("I3").Value = 0.5650593356818092
$ws.Range("O3").Value = 2.515198469554377

$ws.Range("C4").Value = 0.001002524826958506
$ws.Range("D4").Value = 0.001739784821380042
$ws.Range("E4").Value = 0.3412096363495891
$ws.Range("F4").Value = 0.7157517958632837
$ws.Range("G4").Value = 0.002349859887774913
$ws.Range("I4").Value = 0.5371031260879846
$ws.Range("O4").Value = 2.391593294181177

$ws.Range("C5").Value = 0.0009698708151830715
$ws.Range("D5").Value = 0.001693152239104379
$ws.Range("E5").Value = 0.3274337357479027
$ws.Range("F5").Value = 0.7008764035240773
$ws.Range("G5").Value = 0.002351181983142574
$ws.Range("I5").Value = 0.5258055262832073
$ws.Range("O5").Value = 2.341647831043076

$ws.Range("C6").Value = 0.0009644540515694189
$ws.Range("D6").Value = 0.001685434953671461
$ws.Range("E6").Value = 0.3251480654235195
$ws.Range("F6").Value = 0.6984139299659091
$ws.Range("G6").Value = 0.00235140390308072
$ws.Range("I6").Value = 0.5239352555636003
$ws.Range("O6").Value = 2.333379901361866

$ws.Range("C7").Value = 0.001002084080935717
$ws.Range("D7").Value = 0.001739154157165856
$ws.Range("E7").Value = 0.3410237266118372
$ws.Range("F7").Value = 0.7155506723287175
$ws.Range("G7").Value = 0.002349877558024579
$ws.Range("I7").Value = 0.5369503806425797
$ws.Range("O7").Value = 2.390918000611464

$ws.Range("C8").Value = 0.001168915712549179
$ws.Range("D8").Value = 0.001980151755946125
$ws.Range("E8").Value = 0.4113632681818444
$ws.Range("F8").Value = 0.7922838023181242
$ws.Range("G8").Value = 0.002343488075062572
$ws.Range("I8").Value = 0.5952188533536571
$ws.Range("O8").Value = 2.64856492175636

$ws.Range("C9").Value = 0.001498820731356432
$ws.Range("D9").Value = 0.002468573376770422
$ws.Range("E9").Value = 0.5505299427433954
$ws.Range("F9").Value = 0.9468763800419566
$ws.Range("G9").Value = 0.002332191884465699
$ws.Range("I9").Value = 0.7125798299845343
$ws.Range("O9").Value = 3.167691982119493

$ws.Range("C10").Value = 0.001742874980003251
$ws.Range("D10").Value = 0.00283866659964005
$ws.Range("E10").Value = 0.6537366696595797
$ws.Range("F10").Value = 1.063095971564564
$ws.Range("G10").Value = 0.002324635282982213
$ws.Range("I10").Value = 0.800792498283954
$ws.Range("O10").Value = 3.558010637400059

$ws.Range("C11").Value = 0.001854267907294371
$ws.Range("D11").Value = 0.003009832214587504
$ws.Range("E11").Value = 0.7009434488239776
$ws.Range("F11").Value = 1.116569624944702
$ws.Range("G11").Value = 0.002321356869915748
$ws.Range("I11").Value = 0.8413766745340467
$ws.Range("O11").Value = 3.737614285049972

$ws.Range("C12").Value = 0.00189650272443842
$ws.Range("D12").Value = 0.003075081494696263
$ws.Range("E12").Value = 0.7188601408738862
$ws.Range("F12").Value = 1.136907811103981
$ws.Range("G12").Value = 0.002320138148464447
$ws.Range("I12").Value = 0.856812041308217
$ws.Range("O12").Value = 3.805927218431918

$ws.Range("C13").Value = 0.001887404365920275
$ws.Range("D13").Value = 0.003061009190890474
$ws.Range("E13").Value = 0.7149996029872199
$ws.Range("F13").Value = 1.132523635391351
$ws.Range("G13").Value = 0.002320399612671794
$ws.Range("I13").Value = 0.8534847543661783
$ws.Range("O13").Value = 3.791201313206784

$ws.Range("C14").Value = 0.001857741538231039
$ws.Range("D14").Value = 0.003015191475782331
$ws.Range("E14").Value = 0.7024166338524509
$ws.Range("F14").Value = 1.118241067548098
$ws.Range("G14").Value = 0.002321256149805728
$ws.Range("I14").Value = 0.8426451997138287
$ws.Range("O14").Value = 3.743228361005492

$ws.Range("C15").Value = 0.001839579048375839
$ws.Range("D15").Value = 0.002987184000993182
$ws.Range("E15").Value = 0.6947145817116791
$ws.Range("F15").Value = 1.109504213549883
$ws.Range("G15").Value = 0.002321783761155086
$ws.Range("I15").Value = 0.8360144317726252
$ws.Range("O15").Value = 3.71388293500479

$ws.Range("C16").Value = 0.001735602603186948
$ws.Range("D16").Value = 0.002827539475511287
$ws.Range("E16").Value = 0.6506570735473076
$ws.Range("F16").Value = 1.059613674597358
$ws.Range("G16").Value = 0.002324852725136913
$ws.Range("I16").Value = 0.7981495243672896
$ws.Range("O16").Value = 3.546314856958588

$ws.Range("C17").Value = 0.001671911038179985
$ws.Range("D17").Value = 0.002730342243278727
$ws.Range("E17").Value = 0.6236976324255892
$ws.Range("F17").Value = 1.029163690499303
$ws.Range("G17").Value = 0.002326776090929563
$ws.Range("I17").Value = 0.7750384128583931
$ws.Range("O17").Value = 3.444046001339359

$ws.Range("C18").Value = 0.001635312328890137
$ws.Range("D18").Value = 0.002674699540115455
$ws.Range("E18").Value = 0.6082153024573813
$ws.Range("F18").Value = 1.011706473618602
$ws.Range("G18").Value = 0.002327897345802834
$ws.Range("I18").Value = 0.7617883241393031
$ws.Range("O18").Value = 3.385415803844865

$ws.Range("C19").Value = 0.001622926667334212
$ws.Range("D19").Value = 0.002655904110092422
$ws.Range("E19").Value = 0.602977273329472
$ws.Range("F19").Value = 1.005805459256607
$ws.Range("G19").Value = 0.002328279561562133
$ws.Range("I19").Value = 0.7573093800281043
$ws.Range("O19").Value = 3.365597423164161

$ws.Range("C20").Value = 0.001678687496838194
$ws.Range("D20").Value = 0.002740661672977041
$ws.Range("E20").Value = 0.6265649991966171
$ws.Range("F20").Value = 1.032399250097285
$ws.Range("G20").Value = 0.002326569795314328
$ws.Range("I20").Value = 0.7774941886937228
$ws.Range("O20").Value = 3.454912760925026

$ws.Range("C21").Value = 0.001866452803508167
$ws.Range("D21").Value = 0.003028637276567991
$ws.Range("E21").Value = 0.7061114294735376
$ws.Range("F21").Value = 1.122433776554914
$ws.Range("G21").Value = 0.00232100394812837
$ws.Range("I21").Value = 0.8458272088643071
$ws.Range("O21").Value = 3.757310958386086

$ws.Range("C22").Value = 0.001989475412958797
$ws.Range("D22").Value = 0.003219381215675554
$ws.Range("E22").Value = 0.758337397378142
$ws.Range("F22").Value = 1.181795479378252
$ws.Range("G22").Value = 0.00231749884163381
$ws.Range("I22").Value = 0.8908781354258366
$ws.Range("O22").Value = 3.956702928797938

$ws.Range("C23").Value = 0.001923788025923301
$ws.Range("D23").Value = 0.003117335930284071
$ws.Range("E23").Value = 0.7304404936381701
$ws.Range("F23").Value = 1.150064900034494
$ws.Range("G23").Value = 0.002319357505463993
$ws.Range("I23").Value = 0.8667973049203681
$ws.Range("O23").Value = 3.850120629423088

$ws.Range("C24").Value = 0.001675623802700699
$ws.Range("D24").Value = 0.002735995522193946
$ws.Range("E24").Value = 0.6252686099535794
$ws.Range("F24").Value = 1.030936301482171
$ws.Range("G24").Value = 0.002326663013292429
$ws.Range("I24").Value = 0.7763838179350273
$ws.Range("O24").Value = 3.449999383836428

$ws.Range("C25").Value = 0.001409279053683932
$ws.Range("D25").Value = 0.002334584835189446
$ws.Range("E25").Value = 0.5127283134970355
$ws.Range("F25").Value = 0.9046007503700224
$ws.Range("G25").Value = 0.002335116705033783
$ws.Range("I25").Value = 0.6804888268094231
$ws.Range("O25").Value = 3.025721714252654

